# Generate Report for Handoff
#
# The localization pipeline produced a new handoff round for a renamed
# source file (33eab81a-...md -> 2884f280-...md), with fresh xliff
# artifact hashes and updated handoff timestamps. Reflect the new
# report values across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "33eab81a-b54b-47dd-b3ba-274e0ac54933"
$newGuid = "2884f280-caf2-4e65-8a5f-da671b7c46a8"

$oldFileName = "$oldGuid.md"
$newFileName = "$newGuid.md"

$oldPathAndName = "e2e\$oldFileName"
$newPathAndName = "e2e\$newFileName"

$oldHoDate = "2016-09-05 11:24:20"
$newHoDate = "2016-09-05 11:24:59"

$oldZhHash = "8d6f5864d1238f659cb9debd174456916cd23c3b"
$newZhHash = "a25ca0b35341d4bbd284fc6a753823d81cc73fb4"
$oldDeHash = "8d6f5864d1238f659cb9debd174456916cd23c3b"
$newDeHash = "a25ca0b35341d4bbd284fc6a753823d81cc73fb4"

$oldZhHandoffFile = "$oldGuid.$oldZhHash.zh-cn.xlf"
$newZhHandoffFile = "$newGuid.$newZhHash.zh-cn.xlf"
$oldDeHandoffFile = "$oldGuid.$oldDeHash.de-de.xlf"
$newDeHandoffFile = "$newGuid.$newDeHash.de-de.xlf"

$oldZhHandoffDate = "2016-09-05 11:24:13"
$newZhHandoffDate = "2016-09-05 11:24:54"

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathAndName
$wsOverview.Range("G2").Value = $newHoDate

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newFileName
$wsZhCn.Range("G2").Value = $newZhHandoffFile
$wsZhCn.Range("H2").Value = $newZhHandoffDate

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newFileName
$wsDeDe.Range("G2").Value = $newDeHandoffFile
$wsDeDe.Range("H2").Value = $newHoDate
